$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# New column E: flag (1/0) mirroring the existing D-column trigger logic,
# using absolute column references ($B / $C) instead of relative ones.
$ws.Range("E4").Formula = "=IF(AND((`$B4>`$C4),(ABS(`$B4-`$C4)>0.001)), 1,0)"
$ws.Range("E5:E68").Formula = "=IF(AND((`$B5>`$C5),(ABS(`$B5-`$C5)>0.001)), 1,0)"
$ws.Range("E69:E131").Formula = "=IF(AND((`$B69>`$C69),(ABS(`$B69-`$C69)>0.001)), 1,0)"

# Move the selection to the newly-added column.
$ws.Range("E5").Select() | Out-Null
